# Commit: "commit all dao(NOT FINISH)"
# Adds new rows of data to the "Wreath" sheet and to the "Package" sheet,
# and leaves "Package" as the active sheet/tab.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Wreath sheet: add design #1 and #2 rows under the existing header
# ---------------------------------------------------------------
$wreath = $wb.Worksheets.Item("Wreath")

$wreath.Range("A2").Value = "แบบที่1"
$wreath.Range("B2").Value = "jsagxvacgahscvhagscagcga"
$wreath.Range("C2").Value = "C:\Users\User\Desktop\Background\black cat.png"
$wreath.Range("D2").Value = "ผ้า, ดอกไม้"
$wreath.Range("E2").Value = "200/300"
$wreath.Range("F2").Value = "ฟ้า, ม่วง"

$wreath.Range("A3").Value = "แบบที่2"
$wreath.Range("B3").Value = "cnscabchaoica"
$wreath.Range("C3").Value = "C:\Users\User\Desktop\Background\_f2u__sunset_summer___wallpaper_background_by_xxbunnyberryxx_dg55r5k-fullview.jpg"
$wreath.Range("D3").Value = "ผ้า, ดอกไม้"
$wreath.Range("E3").Value = "200/300"
$wreath.Range("F3").Value = "ฟ้า, ม่วง"

# ---------------------------------------------------------------
# Package sheet: add header + 3 package rows
# ---------------------------------------------------------------
$package = $wb.Worksheets.Item("Package")

$package.Range("A1").Value = "แพ็คเกจ"
$package.Range("B1").Value = "รายละเอียด"
$package.Range("C1").Value = "pathรูปภาพ"
$package.Range("D1").Value = "ราคา"

$package.Range("A2").Value = "แพ็คเกจ1"
$package.Range("B2").Value = "ห่ก่หา่ด้หด้หดหด่่ห่ดกสหด้า้หด"
$package.Range("C2").Value = "C:\Users\User\Desktop\Background\black cat.png"
$package.Range("D2").Value = "120000"

$package.Range("A3").Value = "แพ็คเกจ2"
$package.Range("B3").Value = "nakjkjfakhfjkabfkaj"
$package.Range("C3").Value = "C:\Users\User\Desktop\Background\_f2u__sunset_summer___wallpaper_background_by_xxbunnyberryxx_dg55r5k-fullview.jpg"
$package.Range("D3").Value = "21000"

$package.Range("A4").Value = "แพ็คเกจ3"
$package.Range("B4").Value = "หดาหสร่่าาาาาาาาาาาหหสวฟฟสหวสสสสสสสสกสวหสกวหส"
$package.Range("C4").Value = "C:\Users\User\Desktop\Background\_f2u__flower_shop___wallpaper_background_freebie_by_xxbunnyberryxx_dfzmd7c-pre.jpg"
$package.Range("D4").Value = "150000"

# Package becomes the active sheet/tab
$package.Activate()
